$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-06 19:12:42", 0.0008),
    @("2023-12-06 19:13:36", 0.003),
    @("2023-12-06 19:14:51", 0.005000000000000001),
    @("2023-12-06 19:14:56", 0.0006000000000000001)
)

$row = 39
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
